$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Row 8: update CPU Type and Panel Type values ---
$ws1.Range("A8").Value = "MZX254"
$ws1.Range("D8").Value = "FIM"

# --- Rows 8/9: update the "40 V load" values; force text so it keeps "0.000" (not a number) ---
$ws1.Range("J8").Value = "'0.000"
$ws1.Range("J9").Value = "'0.000"

# --- Row 4: fill in the Gallery Type value for the Ethernet Connections block ---
$ws1.Range("B4").ClearFormats() | Out-Null
$ws1.Range("B4").Value = "NGC-601/T1457 OR TC-208"

# --- Update the active selection / scroll position on sheet 1 ---
$ws1.Range("B4").Select() | Out-Null
